$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Splněno?" status for rows 26 and 27 from "ne" to "ano"
$ws.Range("E26").Value = "ano"
$ws.Range("E27").Value = "ano"
